$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New backlog entry goes into the first empty row (row 16), matching the
# existing "Story / Priority / Effort / Validation" table layout.
$ws.Range("A16").Value = "As a developer, I want to create a prototype of the UI."
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = "The prototype is created to make it easier to implement in Android Studio."

# Move the selection the way Excel would leave it after typing the new row.
$ws.Application.Goto($ws.Range("A10"))
$ws.Range("A16:D16").Select()
